$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.082.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.032.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.726"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +18.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.023.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.46%  "

$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.773"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.42%  "

$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("E12").Value = "  -2.07%  "

$ws.Range("E13").Value = "  +15.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.679.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.055.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.51%  "

$ws.Range("E18").Value = "  +1.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.066.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "443.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "104.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.37%  "

$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.64%  "

$ws.Range("E32").Value = "  +4.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.130"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "676.06"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0863"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  -0.97%  "

$ws.Range("E40").Value = "  +5.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.152"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.01%  "

$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("E43").Value = "  +3.89%  "

$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("E46").Value = "  +13.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.53%  "

$ws.Range("E49").Value = "  +2.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.22%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.97%  "
